# Update the rolling 30-day report: shift the date window forward by 5 days
# and refresh the Metrics column with newly generated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{Row=2;  Date="2025-05-24"; Value=92.22734806402124},
    @{Row=3;  Date="2025-05-23"; Value=15.856837896144981},
    @{Row=4;  Date="2025-05-22"; Value=71.0903809450719},
    @{Row=5;  Date="2025-05-21"; Value=68.72211276804812},
    @{Row=6;  Date="2025-05-20"; Value=9.254602841671488},
    @{Row=7;  Date="2025-05-19"; Value=64.18599413793459},
    @{Row=8;  Date="2025-05-18"; Value=93.31987493146613},
    @{Row=9;  Date="2025-05-17"; Value=92.85739095585566},
    @{Row=10; Date="2025-05-16"; Value=12.602063892186177},
    @{Row=11; Date="2025-05-15"; Value=44.669219955222},
    @{Row=12; Date="2025-05-14"; Value=96.63287027506601},
    @{Row=13; Date="2025-05-13"; Value=40.58446468419671},
    @{Row=14; Date="2025-05-12"; Value=10.733536532579013},
    @{Row=15; Date="2025-05-11"; Value=20.48577887214841},
    @{Row=16; Date="2025-05-10"; Value=8.0225586267054},
    @{Row=17; Date="2025-05-09"; Value=57.415730786410236},
    @{Row=18; Date="2025-05-08"; Value=41.158514962576675},
    @{Row=19; Date="2025-05-07"; Value=93.85877590125692},
    @{Row=20; Date="2025-05-06"; Value=24.229812416790786},
    @{Row=21; Date="2025-05-05"; Value=81.01886898305807},
    @{Row=22; Date="2025-05-04"; Value=21.140750912515372},
    @{Row=23; Date="2025-05-03"; Value=57.294525188168734},
    @{Row=24; Date="2025-05-02"; Value=95.22876365628775},
    @{Row=25; Date="2025-05-01"; Value=88.97246215136958},
    @{Row=26; Date="2025-04-30"; Value=35.11854227313118},
    @{Row=27; Date="2025-04-29"; Value=40.04282223732667},
    @{Row=28; Date="2025-04-28"; Value=14.4887007107815},
    @{Row=29; Date="2025-04-27"; Value=46.74362776943893},
    @{Row=30; Date="2025-04-26"; Value=38.63292819991908},
    @{Row=31; Date="2025-04-25"; Value=19.740318025348625}
)

# Style of an existing text cell in column A (General format) used to make
# sure newly written date cells keep the same look (no auto date formatting).
$textStyle = $ws.Cells.Item(1, 1).Style

foreach ($r in $rows) {
    # Prefix with an apostrophe so Excel stores the date as literal text
    # instead of auto-converting it into a date serial number.
    $ws.Cells.Item($r.Row, 1).Value = "'" + $r.Date
    $ws.Cells.Item($r.Row, 1).Style = $textStyle

    $ws.Cells.Item($r.Row, 2).Value = $r.Value
}
